$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 35.88321066666667
$ws.Range("H2").Value = 107.649632
$ws.Range("I2").Value = 0.08317795499144418
$ws.Range("J2").Value = 0.08448843719082051
$ws.Range("M2").Value = 2.084468666666667
$ws.Range("N2").Value = 6.253406
$ws.Range("O2").Value = 0.02757412468035019
$ws.Range("P2").Value = 0.02826332199628731
$ws.Range("Q2").Value = 74.79742829406578
$ws.Range("R2").Value = 673.1768546465921
$ws.Range("S2").Value = 0.002293559301590638
$ws.Range("T2").Value = 0.002387923905287256
# Row 3
$ws.Range("G3").Value = 35.88321066666667
$ws.Range("H3").Value = 107.649632
$ws.Range("I3").Value = 0.08317795499144418
$ws.Range("J3").Value = 0.08448843719082051
$ws.Range("O3").Value = 0.6503186901954997
$ws.Range("P3").Value = 0.6665729829783952
$ws.Range("Q3").Value = 1764.051122640057
$ws.Range("R3").Value = 15876.46010376052
$ws.Range("S3").Value = 0.05409217874317621
$ws.Range("T3").Value = 0.05631770960546801
# Row 4
$ws.Range("G4").Value = 35.88321066666667
$ws.Range("H4").Value = 107.649632
$ws.Range("I4").Value = 0.08317795499144418
$ws.Range("J4").Value = 0.08448843719082051
$ws.Range("M4").Value = 8.752759000000001
$ws.Range("N4").Value = 26.258277
$ws.Range("O4").Value = 0.1157847425689571
$ws.Range("P4").Value = 0.1186787069188703
$ws.Range("Q4").Value = 314.0770951115628
$ws.Range("R4").Value = 2826.693856004065
$ws.Range("S4").Value = 0.009630738106096666
$ws.Range("T4").Value = 0.01002697847540277
# Row 5
$ws.Range("G5").Value = 35.88321066666667
$ws.Range("H5").Value = 107.649632
$ws.Range("I5").Value = 0.08317795499144418
$ws.Range("J5").Value = 0.08448843719082051
$ws.Range("M5").Value = 5.530129000000001
$ws.Range("N5").Value = 11.060258
$ws.Range("O5").Value = 0.07315459761180725
$ws.Range("P5").Value = 0.04998869947289727
$ws.Range("Q5").Value = 198.4387839208427
$ws.Range("R5").Value = 1190.632703525056
$ws.Range("S5").Value = 0.006084849827572112
$ws.Range("T5").Value = 0.004223467095666684
# Row 6
$ws.Range("G6").Value = 35.88321066666667
$ws.Range("H6").Value = 107.649632
$ws.Range("I6").Value = 0.08317795499144418
$ws.Range("J6").Value = 0.08448843719082051
$ws.Range("M6").Value = 10.06683633333333
$ws.Range("N6").Value = 30.200509
$ws.Range("O6").Value = 0.1331678449433857
$ws.Range("P6").Value = 0.1364962886335499
$ws.Range("Q6").Value = 361.2304088958542
$ws.Range("R6").Value = 3251.073680062688
$ws.Range("S6").Value = 0.01107662901300856
$ws.Range("T6").Value = 0.01153235810899579
# Row 7
$ws.Range("I7").Value = 0.03522729558434242
$ws.Range("J7").Value = 0.03578230735158529
$ws.Range("M7").Value = 2.084468666666667
$ws.Range("N7").Value = 6.253406
$ws.Range("O7").Value = 0.02757412468035019
$ws.Range("P7").Value = 0.02826332199628731
$ws.Range("Q7").Value = 31.67799828373688
$ws.Range("R7").Value = 285.101984553632
$ws.Range("S7").Value = 0.0009713618405942076
$ws.Range("T7").Value = 0.001011326874447974
# Row 8
$ws.Range("I8").Value = 0.03522729558434242
$ws.Range("J8").Value = 0.03578230735158529
$ws.Range("O8").Value = 0.6503186901954997
$ws.Range("P8").Value = 0.6665729829783952
$ws.Range("S8").Value = 0.02290896872353927
$ws.Range("T8").Value = 0.02385151934919597
# Row 9
$ws.Range("I9").Value = 0.03522729558434242
$ws.Range("J9").Value = 0.03578230735158529
$ws.Range("M9").Value = 8.752759000000001
$ws.Range("N9").Value = 26.258277
$ws.Range("O9").Value = 0.1157847425689571
$ws.Range("P9").Value = 0.1186787069188703
$ws.Range("Q9").Value = 133.0170556237493
$ws.Range("R9").Value = 1197.153500613744
$ws.Range("S9").Value = 0.004078783350633647
$ws.Range("T9").Value = 0.00424659796705973
# Row 10
$ws.Range("I10").Value = 0.03522729558434242
$ws.Range("J10").Value = 0.03578230735158529
$ws.Range("M10").Value = 5.530129000000001
$ws.Range("N10").Value = 11.060258
$ws.Range("O10").Value = 0.07315459761180725
$ws.Range("P10").Value = 0.04998869947289727
$ws.Range("Q10").Value = 84.04224048662932
$ws.Range("R10").Value = 504.253442919776
$ws.Range("S10").Value = 0.002577038633424764
$ws.Range("T10").Value = 0.00178871100864524
# Row 11
$ws.Range("I11").Value = 0.03522729558434242
$ws.Range("J11").Value = 0.03578230735158529
$ws.Range("M11").Value = 10.06683633333333
$ws.Range("N11").Value = 30.200509
$ws.Range("O11").Value = 0.1331678449433857
$ws.Range("P11").Value = 0.1364962886335499
$ws.Range("Q11").Value = 152.9872956065831
$ws.Range("R11").Value = 1376.885660459248
$ws.Range("S11").Value = 0.004691143036150529
$ws.Range("T11").Value = 0.00488415215223638
# Row 12
$ws.Range("G12").Value = 177.70077
$ws.Range("H12").Value = 533.10231
$ws.Range("I12").Value = 0.4119137160358794
$ws.Range("J12").Value = 0.4184034835782469
$ws.Range("M12").Value = 2.084468666666667
$ws.Range("N12").Value = 6.253406
$ws.Range("O12").Value = 0.02757412468035019
$ws.Range("P12").Value = 0.02826332199628731
$ws.Range("Q12").Value = 370.41168710754
$ws.Range("R12").Value = 3333.70518396786
$ws.Range("S12").Value = 0.0113581601635197
$ws.Range("T12").Value = 0.0118254723807403
# Row 13
$ws.Range("G13").Value = 177.70077
$ws.Range("H13").Value = 533.10231
$ws.Range("I13").Value = 0.4119137160358794
$ws.Range("J13").Value = 0.4184034835782469
$ws.Range("O13").Value = 0.6503186901954997
$ws.Range("P13").Value = 0.6665729829783952
$ws.Range("Q13").Value = 8735.930731630442
$ws.Range("R13").Value = 78623.37658467397
$ws.Range("S13").Value = 0.2678751882860141
$ws.Range("T13").Value = 0.278896458137304
# Row 14
$ws.Range("G14").Value = 177.70077
$ws.Range("H14").Value = 533.10231
$ws.Range("I14").Value = 0.4119137160358794
$ws.Range("J14").Value = 0.4184034835782469
$ws.Range("M14").Value = 8.752759000000001
$ws.Range("N14").Value = 26.258277
$ws.Range("O14").Value = 0.1157847425689571
$ws.Range("P14").Value = 0.1186787069188703
$ws.Range("Q14").Value = 1555.37201392443
$ws.Range("R14").Value = 13998.34812531987
$ws.Range("S14").Value = 0.0476933235718368
$ws.Range("T14").Value = 0.04965558440141712
# Row 15
$ws.Range("G15").Value = 177.70077
$ws.Range("H15").Value = 533.10231
$ws.Range("I15").Value = 0.4119137160358794
$ws.Range("J15").Value = 0.4184034835782469
$ws.Range("M15").Value = 5.530129000000001
$ws.Range("N15").Value = 11.060258
$ws.Range("O15").Value = 0.07315459761180725
$ws.Range("P15").Value = 0.04998869947289727
$ws.Range("Q15").Value = 982.7081814993301
$ws.Range("R15").Value = 5896.249088995981
$ws.Range("S15").Value = 0.03013338214738899
$ws.Range("T15").Value = 0.02091544599900629
# Row 16
$ws.Range("G16").Value = 177.70077
$ws.Range("H16").Value = 533.10231
$ws.Range("I16").Value = 0.4119137160358794
$ws.Range("J16").Value = 0.4184034835782469
$ws.Range("M16").Value = 10.06683633333333
$ws.Range("N16").Value = 30.200509
$ws.Range("O16").Value = 0.1331678449433857
$ws.Range("P16").Value = 0.1364962886335499
$ws.Range("Q16").Value = 1788.88456789731
$ws.Range("R16").Value = 16099.96111107579
$ws.Range("S16").Value = 0.05485366186711981
$ws.Range("T16").Value = 0.05711052265977913
# Row 17
$ws.Range("G17").Value = 20.074196
$ws.Range("H17").Value = 40.148392
$ws.Range("I17").Value = 0.04653236263856699
$ws.Range("J17").Value = 0.0315103250497358
$ws.Range("M17").Value = 2.084468666666667
$ws.Range("N17").Value = 6.253406
$ws.Range("O17").Value = 0.02757412468035019
$ws.Range("P17").Value = 0.02826332199628731
$ws.Range("Q17").Value = 41.84403257052533
$ws.Range("R17").Value = 251.064195423152
$ws.Range("S17").Value = 0.001283089169067115
$ws.Range("T17").Value = 0.0008905864630883607
# Row 18
$ws.Range("G18").Value = 20.074196
$ws.Range("H18").Value = 40.148392
$ws.Range("I18").Value = 0.04653236263856699
$ws.Range("J18").Value = 0.0315103250497358
$ws.Range("O18").Value = 0.6503186901954997
$ws.Range("P18").Value = 0.6665729829783952
$ws.Range("Q18").Value = 986.8656491987788
$ws.Range("R18").Value = 5921.193895192672
$ws.Range("S18").Value = 0.03026086512281489
$ws.Range("T18").Value = 0.02100393136302124
# Row 19
$ws.Range("G19").Value = 20.074196
$ws.Range("H19").Value = 40.148392
$ws.Range("I19").Value = 0.04653236263856699
$ws.Range("J19").Value = 0.0315103250497358
$ws.Range("M19").Value = 8.752759000000001
$ws.Range("N19").Value = 26.258277
$ws.Range("O19").Value = 0.1157847425689571
$ws.Range("P19").Value = 0.1186787069188703
$ws.Range("Q19").Value = 175.704599706764
$ws.Range("R19").Value = 1054.227598240584
$ws.Range("S19").Value = 0.005387737629231838
$ws.Range("T19").Value = 0.003739604631495932
# Row 20
$ws.Range("G20").Value = 20.074196
$ws.Range("H20").Value = 40.148392
$ws.Range("I20").Value = 0.04653236263856699
$ws.Range("J20").Value = 0.0315103250497358
$ws.Range("M20").Value = 5.530129000000001
$ws.Range("N20").Value = 11.060258
$ws.Range("O20").Value = 0.07315459761180725
$ws.Range("P20").Value = 0.04998869947289727
$ws.Range("Q20").Value = 111.012893451284
$ws.Range("R20").Value = 444.0515738051361
$ws.Range("S20").Value = 0.003404056264751062
$ws.Range("T20").Value = 0.00157516016920455
# Row 21
$ws.Range("G21").Value = 20.074196
$ws.Range("H21").Value = 40.148392
$ws.Range("I21").Value = 0.04653236263856699
$ws.Range("J21").Value = 0.0315103250497358
$ws.Range("M21").Value = 10.06683633333333
$ws.Range("N21").Value = 30.200509
$ws.Range("O21").Value = 0.1331678449433857
$ws.Range("P21").Value = 0.1364962886335499
$ws.Range("Q21").Value = 202.0836456552547
$ws.Range("R21").Value = 1212.501873931528
$ws.Range("S21").Value = 0.006196614452702085
$ws.Range("T21").Value = 0.004301042422925714
# Row 22
$ws.Range("G22").Value = 182.547562
$ws.Range("H22").Value = 547.642686
$ws.Range("I22").Value = 0.423148670749767
$ws.Range("J22").Value = 0.4298154468296114
$ws.Range("M22").Value = 2.084468666666667
$ws.Range("N22").Value = 6.253406
$ws.Range("O22").Value = 0.02757412468035019
$ws.Range("P22").Value = 0.02826332199628731
$ws.Range("Q22").Value = 380.5146731653907
$ws.Range("R22").Value = 3424.632058488516
$ws.Range("S22").Value = 0.01166795420557853
$ws.Range("T22").Value = 0.01214801237272341
# Row 23
$ws.Range("G23").Value = 182.547562
$ws.Range("H23").Value = 547.642686
$ws.Range("I23").Value = 0.423148670749767
$ws.Range("J23").Value = 0.4298154468296114
$ws.Range("O23").Value = 0.6503186901954997
$ws.Range("P23").Value = 0.6665729829783952
$ws.Range("Q23").Value = 8974.203414312798
$ws.Range("R23").Value = 80767.83072881518
$ws.Range("S23").Value = 0.2751814893199552
$ws.Range("T23").Value = 0.2865033645234059
# Row 24
$ws.Range("G24").Value = 182.547562
$ws.Range("H24").Value = 547.642686
$ws.Range("I24").Value = 0.423148670749767
$ws.Range("J24").Value = 0.4298154468296114
$ws.Range("M24").Value = 8.752759000000001
$ws.Range("N24").Value = 26.258277
$ws.Range("O24").Value = 0.1157847425689571
$ws.Range("P24").Value = 0.1186787069188703
$ws.Range("Q24").Value = 1597.794816223558
$ws.Range("R24").Value = 14380.15334601202
$ws.Range("S24").Value = 0.04899415991115817
$ws.Range("T24").Value = 0.05100994144349474
# Row 25
$ws.Range("G25").Value = 182.547562
$ws.Range("H25").Value = 547.642686
$ws.Range("I25").Value = 0.423148670749767
$ws.Range("J25").Value = 0.4298154468296114
$ws.Range("M25").Value = 5.530129000000001
$ws.Range("N25").Value = 11.060258
$ws.Range("O25").Value = 0.07315459761180725
$ws.Range("P25").Value = 0.04998869947289727
$ws.Range("Q25").Value = 1009.511566495498
$ws.Range("R25").Value = 6057.069398972989
$ws.Range("S25").Value = 0.03095527073867032
$ws.Range("T25").Value = 0.0214859152003745
# Row 26
$ws.Range("G26").Value = 182.547562
$ws.Range("H26").Value = 547.642686
$ws.Range("I26").Value = 0.423148670749767
$ws.Range("J26").Value = 0.4298154468296114
$ws.Range("M26").Value = 10.06683633333333
$ws.Range("N26").Value = 30.200509
$ws.Range("O26").Value = 0.1331678449433857
$ws.Range("P26").Value = 0.1364962886335499
$ws.Range("Q26").Value = 1837.676429703019
$ws.Range("R26").Value = 16539.08786732717
$ws.Range("S26").Value = 0.05634979657440476
$ws.Range("T26").Value = 0.05866821328961285
